$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a genuinely empty row 2 (placeholder row with no cell data),
# mirroring the existing gaps already present in this sheet's row numbering.
$ws.Rows("2:2").OutlineLevel = 0

# Move the "FSI Agent Governance Framework v1.0 Beta" row from row 12 to row 13.
# Cut() preserves the exact existing cell style (no new style entries created).
$ws.Range("A12").Cut($ws.Range("A13"))
# Remove the leftover cell Cut() leaves behind on the now-empty source cell.
$ws.Range("A12").Clear()

# Re-create empty placeholder rows 11 and 12 (no cell data), like row 2 above.
$ws.Rows("11:11").OutlineLevel = 0
$ws.Rows("12:12").OutlineLevel = 0

# Add the new Control 1.22 row.
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "1.22"
$ws.Range("A10").ClearFormats()
$ws.Range("B10").Value = "Information Barriers for AI Agents"
$ws.Range("C10").Value = "Not Started"
